$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price list refresh (GitHub Actions data pull)

# Row 2
$ws.Range("D2").Value = "57.096.97"
$ws.Range("E2").Value = "  -1.40%  "

# Row 3
$ws.Range("D3").Value = "2.988.35"
$ws.Range("E3").Value = "  -2.28%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.38"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.56"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.12%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.431"
$ws.Range("E8").Value = "  -3.30%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.32"
$ws.Range("E9").Value = "  -4.08%  "

# Row 10
$ws.Range("E10").Value = "  -3.02%  "

# Row 11
$ws.Range("E11").Value = "  -2.15%  "

# Row 12
$ws.Range("D12").Value = "3.491.67"
$ws.Range("E12").Value = "  -2.26%  "

# Row 13
$ws.Range("E13").Value = "  -2.23%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.06"
$ws.Range("E14").Value = "  -1.30%  "

# Row 15
$ws.Range("E15").Value = "  -4.29%  "

# Row 16
$ws.Range("D16").Value = "57.150.09"
$ws.Range("E16").Value = "  -1.31%  "

# Row 17
$ws.Range("E17").Value = "  -1.49%  "

# Row 18
$ws.Range("D18").Value = "2.982.54"
$ws.Range("E18").Value = "  -2.19%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.67"
$ws.Range("E19").Value = "  -2.38%  "

# Row 20
$ws.Range("E20").Value = "  -3.49%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.77"
$ws.Range("E21").Value = "  -5.39%  "

# Row 22
$ws.Range("E22").Value = "  -0.19%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.72"
$ws.Range("E23").Value = "  -0.41%  "

# Row 24
$ws.Range("E24").Value = "  -1.66%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.86"
$ws.Range("E25").Value = "  -2.32%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.165"
$ws.Range("E27").Value = "  -5.65%  "

# Row 28
$ws.Range("E28").Value = "  -6.97%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.58"
$ws.Range("E29").Value = "  -5.20%  "

# Row 30
$ws.Range("E30").Value = "  -2.77%  "

# Row 31
$ws.Range("E31").Value = "  -4.62%  "

# Row 32
$ws.Range("E32").Value = "  -5.82%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.24"
$ws.Range("E33").Value = "  -3.95%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.47"
$ws.Range("E34").Value = "  -1.73%  "

# Row 35
$ws.Range("E35").Value = "  -2.99%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.79"
$ws.Range("E36").Value = "  -2.19%  "

# Row 37
$ws.Range("E37").Value = "  -5.54%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.18"
$ws.Range("E38").Value = "  -5.89%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0667"
$ws.Range("E39").Value = "  -3.71%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.88"

# Row 41
$ws.Range("D41").Value = "3.018.87"
$ws.Range("E41").Value = "  -2.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.27%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.76"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.642"
$ws.Range("E44").Value = "  -2.71%  "

# Row 45
$ws.Range("D45").Value = "2.200.92"
$ws.Range("E45").Value = "  -6.03%  "

# Row 46
$ws.Range("E46").Value = "  -5.89%  "

# Row 47
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.98"
$ws.Range("E47").Value = "  -0.62%  "

# Row 48
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.943"
$ws.Range("E48").Value = "  -8.67%  "

# Row 49
$ws.Range("E49").Value = "  -4.62%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.35"
$ws.Range("E50").Value = "  -3.78%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.83"
$ws.Range("E51").Value = "  -10.22%  "
